# Update the "Point source runs" sheet: replace the computed/formula
# D and E columns (activity results per run) with the newly re-measured
# hard numeric values, refresh the dependent average formulas, and
# highlight the updated overall average (G2) in bold red.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Point source runs")

# New measured values for D (col 4) and E (col 5), rows 2-11.
$newValues = @{
    2  = @(128883, 4743.82)
    3  = @(128782, 4740.12)
    4  = @(128879, 4743.68)
    5  = @(129238, 4756.8900000000003)
    6  = @(128671, 4736.04)
    7  = @(128922, 4745.25)
    8  = @(129029, 4749.1899999999996)
    9  = @(128887, 4743.95)
    10 = @(128787, 4740.3100000000004)
    11 = @(128947, 4746.17)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $dCell = $ws.Range("D$row")
    $eCell = $ws.Range("E$row")

    $dCell.Value = $vals[0]
    $dCell.NumberFormat = "0.00"

    $eCell.Value = $vals[1]
    $eCell.NumberFormat = "0.00"
}

# The overall average of D2:D11 now highlighted in bold red.
$g2 = $ws.Range("G2")
$g2.NumberFormat = "0.00"
$g2.Font.Bold = $true
$g2.Font.Color = 255

# Extend the average of E down through all ten runs (was only E2:E6)
# and give it the same numeric formatting as the rest of the column.
$e14 = $ws.Range("E14")
$e14.Formula = "=AVERAGE(E2:E11)"
$e14.NumberFormat = "0.00"

# Reflect that this sheet was the one last worked in.
$ws.Activate()
$ws.Range("H25").Select()
